# Update recomputed NATMI TPM-based LR-pair statistics (Fgf2-Fgfr1)
# for rows 2-17 in columns E-T, reflecting the new TPM input values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("M2").Value = 7.955277333333332
$ws.Range("N2").Value = 23.865832
$ws.Range("O2").Value = 0.05015625076675284
$ws.Range("P2").Value = 0.05015625076675283
$ws.Range("Q2").Value = 2.491425799976
$ws.Range("R2").Value = 22.422832199784
$ws.Range("S2").Value = 0.001357962732479566
$ws.Range("T2").Value = 0.001357962732479566
$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("M3").Value = 82.48060333333333
$ws.Range("O3").Value = 0.520021823355633
$ws.Range("P3").Value = 0.520021823355633
$ws.Range("Q3").Value = 25.83119287133
$ws.Range("R3").Value = 232.48073584197
$ws.Range("S3").Value = 0.01407940676182124
$ws.Range("T3").Value = 0.01407940676182124
$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 66.90297433333333
$ws.Range("N4").Value = 200.708923
$ws.Range("O4").Value = 0.4218083439585467
$ws.Range("P4").Value = 0.4218083439585465
$ws.Range("Q4").Value = 20.952606598739
$ws.Range("R4").Value = 188.573459388651
$ws.Range("S4").Value = 0.01142031157807995
$ws.Range("T4").Value = 0.01142031157807994
$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("M5").Value = 1.271033333333333
$ws.Range("N5").Value = 3.8131
$ws.Range("O5").Value = 0.008013581919067616
$ws.Range("P5").Value = 0.008013581919067614
$ws.Range("Q5").Value = 0.3980609483
$ws.Range("R5").Value = 3.5825485347
$ws.Range("S5").Value = 0.000216964893376348
$ws.Range("T5").Value = 0.000216964893376348
$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("M6").Value = 7.955277333333332
$ws.Range("N6").Value = 23.865832
$ws.Range("O6").Value = 0.05015625076675284
$ws.Range("P6").Value = 0.05015625076675283
$ws.Range("Q6").Value = 64.25621757828978
$ws.Range("R6").Value = 578.3059582046079
$ws.Range("S6").Value = 0.03502313767572627
$ws.Range("T6").Value = 0.03502313767572627
$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("M7").Value = 82.48060333333333
$ws.Range("O7").Value = 0.520021823355633
$ws.Range("P7").Value = 0.520021823355633
$ws.Range("Q7").Value = 666.2107896060711
$ws.Range("R7").Value = 5995.89710645464
$ws.Range("S7").Value = 0.3631211590847075
$ws.Range("T7").Value = 0.3631211590847075
$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 66.90297433333333
$ws.Range("N8").Value = 200.708923
$ws.Range("O8").Value = 0.4218083439585467
$ws.Range("P8").Value = 0.4218083439585465
$ws.Range("Q8").Value = 540.3874554296791
$ws.Range("R8").Value = 4863.487098867112
$ws.Range("S8").Value = 0.2945405902034232
$ws.Range("T8").Value = 0.2945405902034232
$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("M9").Value = 1.271033333333333
$ws.Range("N9").Value = 3.8131
$ws.Range("O9").Value = 0.008013581919067616
$ws.Range("P9").Value = 0.008013581919067614
$ws.Range("Q9").Value = 10.26636671404445
$ws.Range("R9").Value = 92.39730042640001
$ws.Range("S9").Value = 0.005595728917865167
$ws.Range("T9").Value = 0.005595728917865167
$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("M10").Value = 7.955277333333332
$ws.Range("N10").Value = 23.865832
$ws.Range("O10").Value = 0.05015625076675284
$ws.Range("P10").Value = 0.05015625076675283
$ws.Range("Q10").Value = 23.05237041979822
$ws.Range("R10").Value = 207.471333778184
$ws.Range("S10").Value = 0.01256479720395521
$ws.Range("T10").Value = 0.01256479720395521
$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("M11").Value = 82.48060333333333
$ws.Range("O11").Value = 0.520021823355633
$ws.Range("P11").Value = 0.520021823355633
$ws.Range("Q11").Value = 239.0078108932189
$ws.Range("R11").Value = 2151.07029803897
$ws.Range("S11").Value = 0.1302722721935534
$ws.Range("T11").Value = 0.1302722721935534
$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 66.90297433333333
$ws.Range("N12").Value = 200.708923
$ws.Range("O12").Value = 0.4218083439585467
$ws.Range("P12").Value = 0.4218083439585465
$ws.Range("Q12").Value = 193.8678039615279
$ws.Range("R12").Value = 1744.810235653751
$ws.Range("S12").Value = 0.1056685103003852
$ws.Range("T12").Value = 0.1056685103003852
$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("M13").Value = 1.271033333333333
$ws.Range("N13").Value = 3.8131
$ws.Range("O13").Value = 0.008013581919067616
$ws.Range("P13").Value = 0.008013581919067614
$ws.Range("Q13").Value = 3.683131333855556
$ws.Range("R13").Value = 33.1481820047
$ws.Range("S13").Value = 0.002007507143199601
$ws.Range("T13").Value = 0.002007507143199601
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("M14").Value = 7.955277333333332
$ws.Range("N14").Value = 23.865832
$ws.Range("O14").Value = 0.05015625076675284
$ws.Range("P14").Value = 0.05015625076675283
$ws.Range("Q14").Value = 2.220609597235555
$ws.Range("R14").Value = 19.98548637512
$ws.Range("S14").Value = 0.001210353154591797
$ws.Range("T14").Value = 0.001210353154591797
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("M15").Value = 82.48060333333333
$ws.Range("O15").Value = 0.520021823355633
$ws.Range("P15").Value = 0.520021823355633
$ws.Range("Q15").Value = 23.02336067912222
$ws.Range("R15").Value = 207.2102461121
$ws.Range("S15").Value = 0.01254898531555087
$ws.Range("T15").Value = 0.01254898531555088
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 66.90297433333333
$ws.Range("N16").Value = 200.708923
$ws.Range("O16").Value = 0.4218083439585467
$ws.Range("P16").Value = 0.4218083439585465
$ws.Range("Q16").Value = 18.67507324549222
$ws.Range("R16").Value = 168.07565920943
$ws.Range("S16").Value = 0.01017893187665832
$ws.Range("T16").Value = 0.01017893187665832
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("M17").Value = 1.271033333333333
$ws.Range("N17").Value = 3.8131
$ws.Range("O17").Value = 0.008013581919067616
$ws.Range("P17").Value = 0.008013581919067614
$ws.Range("Q17").Value = 0.3547920078888889
$ws.Range("R17").Value = 3.193128071
$ws.Range("S17").Value = 0.0001933809646264996
$ws.Range("T17").Value = 0.0001933809646264996
